$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Power Hour: Python June 21, 2020"
#    -> "Power Hour: Python June" | " Aug 14th" | [_GoBack bookmark] | ", 2020"
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(2)
$titleStart = $titlePara.Range.Start

# Replace "21" with "Aug 14th" (keeps surrounding spaces intact, single run)
$replaceRange = $d.Range($titleStart + 24, $titleStart + 26)
$replaceRange.Text = "Aug 14th"

# Re-fetch the paragraph start (unchanged, but stay consistent)
$titlePara = $d.Paragraphs.Item(2)
$titleStart = $titlePara.Range.Start

# Move the existing (single) "_GoBack" bookmark so it sits right after
# " Aug 14th" and before ", 2020". Adding it first forces the run split
# at that boundary without leaving a residual xml:space artifact on the
# trailing run.
$finalBmRange = $d.Range($titleStart + 32, $titleStart + 32)
$d.Bookmarks.Add("_GoBack", $finalBmRange)

# Now split off "Power Hour: Python June" from " Aug 14th" using a
# temporary bookmark (creates the run break, then remove the temp mark).
$titlePara = $d.Paragraphs.Item(2)
$titleStart = $titlePara.Range.Start
$tempBmRange = $d.Range($titleStart + 23, $titleStart + 23)
$d.Bookmarks.Add("ZZTempSplit", $tempBmRange)
$d.Bookmarks.Item("ZZTempSplit").Delete()

Write-Output "title updated"
